$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 564 (pushes the existing rows 564-665 down to 565-666)
$ws.Rows.Item(564).Insert()

# Populate the newly inserted row with the weekly Perejil price record
$ws.Cells.Item(564, 1).Value = 6
$ws.Cells.Item(564, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(564, 3).Value = "Metropolitana"
$ws.Cells.Item(564, 4).Value = 44951
$ws.Cells.Item(564, 5).Value = 13
$ws.Cells.Item(564, 6).Value = 100112044
$ws.Cells.Item(564, 7).Value = "Perejil"
$ws.Cells.Item(564, 8).Value = "Sin especificar"
$ws.Cells.Item(564, 9).Value = "Primera"
$ws.Cells.Item(564, 10).Value = 190
$ws.Cells.Item(564, 11).Value = 17000
$ws.Cells.Item(564, 12).Value = 18000
$ws.Cells.Item(564, 13).Value = 17421
$ws.Cells.Item(564, 14).Value = "`$/docena de atados"
$ws.Cells.Item(564, 15).Value = "Región Metropolitana"
$ws.Cells.Item(564, 16).Value = 5807
$ws.Cells.Item(564, 17).Value = 3
$ws.Cells.Item(564, 18).Value = "Hortaliza"
